# Insert a new weekly price record at the top of the Espinaca data block.
# The existing rows 372-399 shift down to 373-400, and the new row 372
# holds the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(372).Insert()

$ws.Range("A372").Value = 8
$ws.Range("B372").Value = "Terminal La Palmera de La Serena"
$ws.Range("C372").Value = "Coquimbo"
$ws.Range("D372").Value = 45021
$ws.Range("E372").Value = 4
$ws.Range("F372").Value = 100112012
$ws.Range("G372").Value = "Espinaca"
$ws.Range("H372").Value = "Sin especificar"
$ws.Range("I372").Value = "Primera"
$ws.Range("J372").Value = 2000
$ws.Range("K372").Value = 450
$ws.Range("L372").Value = 500
$ws.Range("M372").Value = 475
$ws.Range("N372").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O372").Value = "Provincia del Elquí"
$ws.Range("P372").Value = 950
$ws.Range("Q372").Value = 0.5
$ws.Range("R372").Value = "Hortaliza"
